$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3100.9524
$ws.Range("I64").Value = 3121.2144
$ws.Range("K64").Value = 3121.2144
$ws.Range("M64").Value = -2873.2144
$ws.Range("H67").Value = 3100.9524
$ws.Range("I67").Value = 3121.2144
$ws.Range("K67").Value = 3121.2144
$ws.Range("M67").Value = -2263.2144
$ws.Range("H74").Value = 3350.7368
$ws.Range("I74").Value = 3222
$ws.Range("J74").Value = 3466.6
$ws.Range("K74").Value = 3222
$ws.Range("L74").Value = 3466.6
$ws.Range("M74").Value = -2286
$ws.Range("N74").Value = -5338.6
$ws.Range("H77").Value = 3350.7368
$ws.Range("I77").Value = 3222
$ws.Range("J77").Value = 3466.6
$ws.Range("K77").Value = 16110
$ws.Range("L77").Value = 17333
$ws.Range("M77").Value = -11430
$ws.Range("N77").Value = -26693
$ws.Range("H103").Value = 3756616
$ws.Range("I103").Value = 6677995
$ws.Range("J103").Value = 557.1429000000001
$ws.Range("K103").Value = 20033985
$ws.Range("L103").Value = 1671.4287
$ws.Range("M103").Value = -20033399
$ws.Range("N103").Value = -2843.4287
$ws.Range("H112").Value = 1357.8269
$ws.Range("J112").Value = 1215.7709
$ws.Range("L112").Value = 3647.3127
$ws.Range("N112").Value = -5863.3127
$ws.Range("H132").Value = 5689051
$ws.Range("I132").Value = 6496.643
$ws.Range("J132").Value = 15633522
$ws.Range("K132").Value = 19489.929
$ws.Range("L132").Value = 46900566
$ws.Range("M132").Value = -16959.929
$ws.Range("N132").Value = -46905626
$ws.Range("H135").Value = 19232744
$ws.Range("I135").Value = 1968.9546
$ws.Range("J135").Value = 125002010
$ws.Range("K135").Value = 17720.5914
$ws.Range("L135").Value = 1125018090
$ws.Range("M135").Value = -15185.5914
$ws.Range("N135").Value = -1125023160
$ws.Range("H138").Value = 1705.1735
$ws.Range("I138").Value = 1058.2307
$ws.Range("J138").Value = 2132.8135
$ws.Range("K138").Value = 3174.6921
$ws.Range("L138").Value = 6398.440500000001
$ws.Range("M138").Value = 1965.3079
$ws.Range("N138").Value = -16678.4405

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8655.879000000001
$ws.Range("I32").Value = 9507.861999999999
$ws.Range("K32").Value = 9507.861999999999
$ws.Range("M32").Value = -9220.861999999999
$ws.Range("H61").Value = 9617322
$ws.Range("I61").Value = 11629659
$ws.Range("J61").Value = 2823.7778
$ws.Range("K61").Value = 11629659
$ws.Range("L61").Value = 2823.7778
$ws.Range("M61").Value = -11629447
$ws.Range("N61").Value = -3247.7778
$ws.Range("H97").Value = 4165.483
$ws.Range("I97").Value = 5755.8
$ws.Range("J97").Value = 631.44446
$ws.Range("K97").Value = 5755.8
$ws.Range("L97").Value = 631.44446
$ws.Range("M97").Value = -5259.8
$ws.Range("N97").Value = -1623.44446
$ws.Range("H123").Value = 31804
$ws.Range("J123").Value = 31804
$ws.Range("L123").Value = 31804
$ws.Range("N123").Value = -41604
$ws.Range("H132").Value = 5954482.5
$ws.Range("I132").Value = 10871238
$ws.Range("K132").Value = 32613714
$ws.Range("M132").Value = -32611184
$ws.Range("H136").Value = 9617322
$ws.Range("I136").Value = 11629659
$ws.Range("J136").Value = 2823.7778
$ws.Range("K136").Value = 34888977
$ws.Range("L136").Value = 8471.3334
$ws.Range("M136").Value = -34886427
$ws.Range("N136").Value = -13571.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 570.41174
$ws.Range("I94").Value = 460.5
$ws.Range("J94").Value = 727.4286
$ws.Range("K94").Value = 460.5
$ws.Range("L94").Value = 727.4286
$ws.Range("M94").Value = -9.5
$ws.Range("N94").Value = -1629.4286
$ws.Range("H134").Value = 5292.082
$ws.Range("I134").Value = 4567.525
$ws.Range("J134").Value = 6672.1904
$ws.Range("K134").Value = 13702.575
$ws.Range("L134").Value = 20016.5712
$ws.Range("M134").Value = -11167.575
$ws.Range("N134").Value = -25086.5712
$ws.Range("H135").Value = 38853.625
$ws.Range("J135").Value = 38853.625
$ws.Range("L135").Value = 38853.625
$ws.Range("N135").Value = -48993.625
$ws.Range("H138").Value = 44949.832
$ws.Range("J138").Value = 44949.832
$ws.Range("L138").Value = 44949.832
$ws.Range("N138").Value = -55229.832

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5750781.5
$ws.Range("I31").Value = 5409.8438
$ws.Range("J31").Value = 12822008
$ws.Range("K31").Value = 5409.8438
$ws.Range("L31").Value = 12822008
$ws.Range("M31").Value = -5114.8438
$ws.Range("N31").Value = -12822598
$ws.Range("H34").Value = 5750781.5
$ws.Range("I34").Value = 5409.8438
$ws.Range("J34").Value = 12822008
$ws.Range("K34").Value = 5409.8438
$ws.Range("L34").Value = 12822008
$ws.Range("M34").Value = -5207.8438
$ws.Range("N34").Value = -12822412

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 799.125
$ws.Range("I5").Value = 516.38464
$ws.Range("J5").Value = 1133.2727
$ws.Range("K5").Value = 1549.15392
$ws.Range("L5").Value = 3399.8181
$ws.Range("M5").Value = -1437.15392
$ws.Range("N5").Value = -3623.8181
$ws.Range("H50").Value = 192.83333
$ws.Range("I50").Value = 235
$ws.Range("J50").Value = 150.66667
$ws.Range("K50").Value = 705
$ws.Range("L50").Value = 452.00001
$ws.Range("M50").Value = -224
$ws.Range("N50").Value = -1414.00001
$ws.Range("H53").Value = 192.83333
$ws.Range("I53").Value = 235
$ws.Range("J53").Value = 150.66667
$ws.Range("K53").Value = 705
$ws.Range("L53").Value = 452.00001
$ws.Range("M53").Value = -224
$ws.Range("N53").Value = -1414.00001
$ws.Range("H103").Value = 3973.625
$ws.Range("I103").Value = 1770
$ws.Range("J103").Value = 4414.35
$ws.Range("K103").Value = 5310
$ws.Range("L103").Value = 13243.05
$ws.Range("M103").Value = -4431
$ws.Range("N103").Value = -15001.05
$ws.Range("H114").Value = 1221.1333
$ws.Range("I114").Value = 320.6
$ws.Range("J114").Value = 3022.2
$ws.Range("K114").Value = 961.8000000000001
$ws.Range("L114").Value = 9066.599999999999
$ws.Range("M114").Value = 2292.2
$ws.Range("N114").Value = -15574.6
$ws.Range("H118").Value = 1998.1765
$ws.Range("I118").Value = 329
$ws.Range("J118").Value = 2102.5
$ws.Range("K118").Value = 987
$ws.Range("L118").Value = 6307.5
$ws.Range("M118").Value = 256
$ws.Range("N118").Value = -8793.5
$ws.Range("H135").Value = 799.125
$ws.Range("I135").Value = 516.38464
$ws.Range("J135").Value = 1133.2727
$ws.Range("K135").Value = 4647.46176
$ws.Range("L135").Value = 10199.4543
$ws.Range("M135").Value = -2112.46176
$ws.Range("N135").Value = -15269.4543

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1233.3334
$ws.Range("I113").Value = 1450
$ws.Range("J113").Value = 800
$ws.Range("K113").Value = 1450
$ws.Range("L113").Value = 800
$ws.Range("M113").Value = 720
$ws.Range("N113").Value = -5140
$ws.Range("H132").Value = 4585.5454
$ws.Range("I132").Value = 3579.8667
$ws.Range("K132").Value = 10739.6001
$ws.Range("M132").Value = -8209.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7999.9414
$ws.Range("I7").Value = 12839.8
$ws.Range("J7").Value = 5983.3335
$ws.Range("K7").Value = 12839.8
$ws.Range("L7").Value = 5983.3335
$ws.Range("M7").Value = -12727.8
$ws.Range("N7").Value = -6207.3335
$ws.Range("H40").Value = 5414.769
$ws.Range("I40").Value = 7091.4165
$ws.Range("J40").Value = 3977.6428
$ws.Range("K40").Value = 7091.4165
$ws.Range("L40").Value = 3977.6428
$ws.Range("M40").Value = -6955.4165
$ws.Range("N40").Value = -4249.6428
$ws.Range("H68").Value = 1600.125
$ws.Range("I68").Value = 1093.6666
$ws.Range("J68").Value = 1904
$ws.Range("K68").Value = 1093.6666
$ws.Range("L68").Value = 1904
$ws.Range("M68").Value = -344.6666
$ws.Range("N68").Value = -3402
$ws.Range("H71").Value = 1600.125
$ws.Range("I71").Value = 1093.6666
$ws.Range("J71").Value = 1904
$ws.Range("K71").Value = 5468.333000000001
$ws.Range("L71").Value = 9520
$ws.Range("M71").Value = -1724.333000000001
$ws.Range("N71").Value = -17008
$ws.Range("H100").Value = 1961.4
$ws.Range("I100").Value = 1701.5
$ws.Range("J100").Value = 2134.6667
$ws.Range("K100").Value = 1701.5
$ws.Range("L100").Value = 2134.6667
$ws.Range("M100").Value = -1160.5
$ws.Range("N100").Value = -3216.6667
$ws.Range("H126").Value = 7999.9414
$ws.Range("I126").Value = 12839.8
$ws.Range("J126").Value = 5983.3335
$ws.Range("K126").Value = 38519.39999999999
$ws.Range("L126").Value = 17950.0005
$ws.Range("M126").Value = -36049.39999999999
$ws.Range("N126").Value = -22890.0005
$ws.Range("H132").Value = 11119102
$ws.Range("I132").Value = 5396.185
$ws.Range("J132").Value = 27789662
$ws.Range("K132").Value = 16188.555
$ws.Range("L132").Value = 83368986
$ws.Range("M132").Value = -13658.555
$ws.Range("N132").Value = -83374046
$ws.Range("H133").Value = 52231
$ws.Range("J133").Value = 52231
$ws.Range("L133").Value = 52231
$ws.Range("N133").Value = -57291

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3430.7693
$ws.Range("I96").Value = 2333.3333
$ws.Range("J96").Value = 3760
$ws.Range("K96").Value = 2333.3333
$ws.Range("L96").Value = 3760
$ws.Range("M96").Value = -960.3332999999998
$ws.Range("N96").Value = -6506
$ws.Range("H100").Value = 1742.7894
$ws.Range("I100").Value = 2093.9092
$ws.Range("J100").Value = 1260
$ws.Range("K100").Value = 4187.8184
$ws.Range("L100").Value = 2520
$ws.Range("M100").Value = -3646.8184
$ws.Range("N100").Value = -3602
$ws.Range("H107").Value = 832.1053000000001
$ws.Range("I107").Value = 1232.4546
$ws.Range("J107").Value = 281.625
$ws.Range("K107").Value = 3697.3638
$ws.Range("L107").Value = 844.875
$ws.Range("M107").Value = -1777.3638
$ws.Range("N107").Value = -4684.875
$ws.Range("H113").Value = 2272.476
$ws.Range("I113").Value = 2486.2222
$ws.Range("J113").Value = 2112.1667
$ws.Range("K113").Value = 7458.6666
$ws.Range("L113").Value = 6336.500100000001
$ws.Range("M113").Value = -5288.6666
$ws.Range("N113").Value = -10676.5001
$ws.Range("H123").Value = 38348.637
$ws.Range("J123").Value = 38348.637
$ws.Range("L123").Value = 38348.637
$ws.Range("N123").Value = -48148.637
$ws.Range("H135").Value = 46715
$ws.Range("J135").Value = 46715
$ws.Range("L135").Value = 46715
$ws.Range("N135").Value = -56855
$ws.Range("H138").Value = 60429
$ws.Range("J138").Value = 60429
$ws.Range("L138").Value = 60429
$ws.Range("N138").Value = -70709
